$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 595
$ws.Range("F4").Value = 1274
$ws.Range("F5").Value = 1125
$ws.Range("F6").Value = 14159
$ws.Range("F7").Value = 15866
$ws.Range("F9").Value = 64
$ws.Range("F11").Value = 196
$ws.Range("F17").Value = 33
$ws.Range("F18").Value = 85
$ws.Range("F20").Value = 1228
$ws.Range("F21").Value = 130
$ws.Range("F23").Value = 20
$ws.Range("F24").Value = 6273
$ws.Range("F26").Value = 1101
$ws.Range("F27").Value = 5605
$ws.Range("F28").Value = 80
$ws.Range("F29").Value = 143
$ws.Range("F30").Value = 132
$ws.Range("F31").Value = 4589
$ws.Range("F32").Value = 7

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 595
$ws.Range("F4").Value = 1274
$ws.Range("F5").Value = 1125
$ws.Range("F6").Value = 14159
$ws.Range("F7").Value = 15866
$ws.Range("F9").Value = 64
$ws.Range("F11").Value = 196
$ws.Range("F17").Value = 33
$ws.Range("F18").Value = 85
$ws.Range("F20").Value = 1228
$ws.Range("F21").Value = 130
$ws.Range("F24").Value = 20
$ws.Range("F25").Value = 6273
$ws.Range("F27").Value = 1101
$ws.Range("F28").Value = 5605
$ws.Range("F29").Value = 80
$ws.Range("F30").Value = 143
$ws.Range("F31").Value = 132
$ws.Range("F32").Value = 4589
$ws.Range("F33").Value = 7
